$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.017.41"
$ws.Range("E2").Value = "  +0.14%  "

$ws.Range("D3").Value = "2.954.08"
$ws.Range("E3").Value = "  +0.90%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "380.99"
$ws.Range("E5").Value = "  +0.87%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.37"
$ws.Range("E6").Value = "  +0.44%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.542"
$ws.Range("E7").Value = "  +1.57%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("E9").Value = "  +0.54%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.36"
$ws.Range("E10").Value = "  -0.06%  "

$ws.Range("E11").Value = "  -0.43%  "

$ws.Range("E12").Value = "  +1.95%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.41"
$ws.Range("E13").Value = "  +3.02%  "

$ws.Range("D14").Value = "3.418.14"
$ws.Range("E14").Value = "  +1.08%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "12.36"
$ws.Range("E15").Value = "  +74.27%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.74"
$ws.Range("E16").Value = "  +5.68%  "

$ws.Range("D17").Value = "2.951.02"
$ws.Range("E17").Value = "  +1.07%  "

$ws.Range("E18").Value = "  +4.28%  "

$ws.Range("D19").Value = "51.035.42"
$ws.Range("E19").Value = "  +0.38%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.08"
$ws.Range("E20").Value = "  -2.34%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.35"
$ws.Range("E21").Value = "  -0.89%  "

$ws.Range("D22").Value = "0.0₃0954"
$ws.Range("E22").Value = "  +0.91%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.39"
$ws.Range("E23").Value = "  +18.63%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.64"
$ws.Range("E24").Value = "  +2.37%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "266.64"
$ws.Range("E25").Value = "  +2.04%  "

$ws.Range("E26").Value = "  -1.54%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "25.81"
$ws.Range("E28").Value = "  +1.24%  "

$ws.Range("E29").Value = "  -1.08%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.93"
$ws.Range("E30").Value = "  -7.27%  "

$ws.Range("E31").Value = "  -4.32%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.48"
$ws.Range("E32").Value = "  +7.43%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "50.68"
$ws.Range("E33").Value = "  -0.02%  "

$ws.Range("E34").Value = "  +2.10%  "

$ws.Range("E35").Value = "  +0.48%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0433"
$ws.Range("E36").Value = "  -3.86%  "

$ws.Range("E37").Value = "  +0.01%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.19"
$ws.Range("E38").Value = "  +7.76%  "

$ws.Range("E41").Value = "  +2.84%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.48"
$ws.Range("E42").Value = "  -3.53%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "119.82"
$ws.Range("E43").Value = "  -0.56%  "

$ws.Range("E44").Value = "  +10.93%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.47"
$ws.Range("E45").Value = "  +1.66%  "

$ws.Range("E46").Value = "  -1.04%  "

$ws.Range("D47").Value = "2.025.07"
$ws.Range("E47").Value = "  +1.29%  "

$ws.Range("E48").Value = "  -1.83%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.258"
$ws.Range("E49").Value = "  -4.64%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0319"
$ws.Range("E50").Value = "  -7.48%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.36"
$ws.Range("E51").Value = "  +7.51%  "

# Row 39/40: Stellar/Celestia swap
$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.70"
$ws.Range("E39").Value = "  +2.72%  "

$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.116"
$ws.Range("E40").Value = "  +1.88%  "
